# Weekly cryptos-list refresh (GitHub Actions scraper run, Sat Jun 24 16:45:52 UTC 2023).
# Updates the Price (D) and Volume(1h) (E) columns for each coin, and fixes the
# ranking order of Quant / PaxDollar (rows 42-43 swap position).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E hold plain-text values (e.g. "16.60", "0.000007770", "  -1.52%  ") in the
# source file. Pre-format the range as Text so Excel does not silently reinterpret
# them as numbers/percentages and strip significant trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.483.04'
$ws.Range('E2').Value = '  -1.52%  '
$ws.Range('D3').Value = '1.877.10'
$ws.Range('E3').Value = '  -1.84%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '232.82'
$ws.Range('E5').Value = '  -5.35%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.4855'
$ws.Range('E7').Value = '  -2.06%  '
$ws.Range('D8').Value = '0.2873'
$ws.Range('E8').Value = '  -4.22%  '
$ws.Range('D9').Value = '0.06617'
$ws.Range('E9').Value = '  -2.39%  '
$ws.Range('D10').Value = '1.878.82'
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('D11').Value = '16.60'
$ws.Range('E11').Value = '  -2.34%  '
$ws.Range('D12').Value = '0.07214'
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('D13').Value = '88.08'
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('D14').Value = '4.951'
$ws.Range('E14').Value = '  -2.39%  '
$ws.Range('D15').Value = '0.6577'
$ws.Range('E15').Value = '  -3.63%  '
$ws.Range('D16').Value = '30.457.40'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '0.000007770'
$ws.Range('E18').Value = '  -3.15%  '
$ws.Range('D19').Value = '12.83'
$ws.Range('E19').Value = '  -2.72%  '
$ws.Range('D20').Value = '2.121.41'
$ws.Range('E20').Value = '  -1.46%  '
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').Value = '4.707'
$ws.Range('E22').Value = '  -3.37%  '
$ws.Range('D23').Value = '185.50'
$ws.Range('E23').Value = '  +5.69%  '
$ws.Range('D24').Value = '5.985'
$ws.Range('E24').Value = '  -1.10%  '
$ws.Range('D25').Value = '9.191'
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('D26').Value = '155.14'
$ws.Range('E26').Value = '  +1.70%  '
$ws.Range('D27').Value = '18.29'
$ws.Range('E27').Value = '  +1.40%  '
$ws.Range('D28').Value = '1.826'
$ws.Range('E28').Value = '  -6.38%  '
$ws.Range('D29').Value = '1.396'
$ws.Range('E29').Value = '  -1.54%  '
$ws.Range('D30').Value = '4.219'
$ws.Range('E30').Value = '  -2.49%  '
$ws.Range('D31').Value = '0.08963'
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('D32').Value = '3.892'
$ws.Range('E32').Value = '  -4.47%  '
$ws.Range('D33').Value = '0.05170'
$ws.Range('E33').Value = '  -2.28%  '
$ws.Range('D34').Value = '0.7297'
$ws.Range('E34').Value = '  -2.33%  '
$ws.Range('D35').Value = '1.071'
$ws.Range('E35').Value = '  -6.24%  '
$ws.Range('D36').Value = '2.697'
$ws.Range('E36').Value = '  +2.04%  '
$ws.Range('D37').Value = '0.01806'
$ws.Range('E37').Value = '  -5.48%  '
$ws.Range('D38').Value = '2.647'
$ws.Range('E38').Value = '  -2.80%  '
$ws.Range('E39').Value = '  -2.75%  '
$ws.Range('D40').Value = '2.007'
$ws.Range('E40').Value = '  -9.24%  '
$ws.Range('D41').Value = '0.4283'
$ws.Range('E41').Value = '  -2.38%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '0.9957'
$ws.Range('E42').Value = '  -0.55%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '103.26'
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('D44').Value = '5.588'
$ws.Range('E44').Value = '  -6.69%  '
$ws.Range('D45').Value = '0.1326'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').Value = '7.178'
$ws.Range('E46').Value = '  -8.35%  '
$ws.Range('D47').Value = '0.05804'
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('E48').Value = '  +0.99%  '
$ws.Range('D49').Value = '1.393'
$ws.Range('E49').Value = '  +0.98%  '
$ws.Range('D50').Value = '0.3852'
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('D51').Value = '32.98'
$ws.Range('E51').Value = '  -1.25%  '
